# Weekly update for "Poroto verde" (Vega Central Mapocho de Santiago)
# Two new price records are inserted right after the existing row 226
# (i.e. they become the new rows 227-228), pushing every following
# record down by two rows. The sheet's dimension grows from A1:R260 to
# A1:R262 automatically once the data is in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows at position 227 (old row 227 and everything
# below shifts down by two, ending at row 262).
$ws.Rows.Item(227).Insert()
$ws.Rows.Item(227).Insert()

# ---- New row 227 ----
$ws.Cells.Item(227, 1).Value  = 9
$ws.Cells.Item(227, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(227, 3).Value  = "Metropolitana"
$ws.Cells.Item(227, 4).Value  = 44476
$ws.Cells.Item(227, 5).Value  = 13
$ws.Cells.Item(227, 6).Value  = 100112031
$ws.Cells.Item(227, 7).Value  = "Poroto verde"
$ws.Cells.Item(227, 8).Value  = "Magnum"
$ws.Cells.Item(227, 9).Value  = "Primera"
$ws.Cells.Item(227, 10).Value = 52
$ws.Cells.Item(227, 11).Value = 43000
$ws.Cells.Item(227, 12).Value = 45000
$ws.Cells.Item(227, 13).Value = 44000
$ws.Cells.Item(227, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(227, 15).Value = "Perú"
$ws.Cells.Item(227, 16).Value = 1760
$ws.Cells.Item(227, 17).Value = 25
$ws.Cells.Item(227, 18).Value = "Hortaliza"

# ---- New row 228 ----
$ws.Cells.Item(228, 1).Value  = 9
$ws.Cells.Item(228, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(228, 3).Value  = "Metropolitana"
$ws.Cells.Item(228, 4).Value  = 44476
$ws.Cells.Item(228, 5).Value  = 13
$ws.Cells.Item(228, 6).Value  = 100112031
$ws.Cells.Item(228, 7).Value  = "Poroto verde"
$ws.Cells.Item(228, 8).Value  = "Sin especificar"
$ws.Cells.Item(228, 9).Value  = "Primera"
$ws.Cells.Item(228, 10).Value = 34
$ws.Cells.Item(228, 11).Value = 37000
$ws.Cells.Item(228, 12).Value = 39000
$ws.Cells.Item(228, 13).Value = 38000
$ws.Cells.Item(228, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(228, 15).Value = "Perú"
$ws.Cells.Item(228, 16).Value = 1520
$ws.Cells.Item(228, 17).Value = 25
$ws.Cells.Item(228, 18).Value = "Hortaliza"

# Keep the date column's display format ("YYYY-MM-DD HH:MM:SS") consistent
# with the rest of column D for the two freshly inserted rows.
$ws.Range("D227:D228").NumberFormat = $ws.Range("D229").NumberFormat
